$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "VALOR MORA" total figure
$ws.Range("E11").Value = 123641

# Update worker/period counters
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 1

# Update Valor Mora for NILSON VILLADIEGO GUERRERO row
$ws.Range("G16").Value = 1725042

# Remove the worker row for 73104559 / LEVER M TORRES VILLADIEGO (row 17),
# shifting the rows below it up by one
$ws.Rows("17:17").Delete()

# Update Valor Mora for JORGE MIGUEL PEREZ MONTERROZA, now on row 17
# after the deletion above
$ws.Range("G17").Value = 1366000
